$d = $word.ActiveDocument

$replacements = @(
    @("711÷6=118, 3", "575÷5=115, 0"),
    @("313÷5=62, 3", "442÷3=147, 1"),
    @("395÷4=98, 3", "869÷6=144, 5"),
    @("308÷9=34, 2", "491÷3=163, 2"),
    @("327÷7=46, 5", "908÷2=454, 0"),
    @("950÷9=105, 5", "679÷4=169, 3"),
    @("627÷2=313, 1", "148÷9=16, 4"),
    @("536÷7=76, 4", "180÷8=22, 4"),
    @("748÷7=106, 6", "627÷7=89, 4"),
    @("371÷8=46, 3", "227÷4=56, 3"),
    @("568÷7=81, 1", "458÷2=229, 0"),
    @("708÷2=354, 0", "894÷9=99, 3"),
    @("439÷5=87, 4", "915÷5=183, 0"),
    @("449÷2=224, 1", "887÷9=98, 5"),
    @("407÷3=135, 2", "675÷5=135, 0"),
    @("250÷4=62, 2", "782÷3=260, 2"),
    @("842÷3=280, 2", "270÷5=54, 0"),
    @("823÷8=102, 7", "759÷4=189, 3"),
    @("543÷2=271, 1", "742÷5=148, 2"),
    @("439÷2=219, 1", "714÷8=89, 2"),
    @("768÷8=96, 0", "373÷2=186, 1"),
    @("427÷8=53, 3", "885÷3=295, 0"),
    @("280÷8=35, 0", "124÷6=20, 4"),
    @("455÷6=75, 5", "605÷3=201, 2"),
    @("944÷8=118, 0", "168÷7=24, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done"
